$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6139386892318726
$ws.Range("B1").Value = 1.138901114463806
$ws.Range("C1").Value = 5.780870914459229
$ws.Range("D1").Value = 1.734243154525757
$ws.Range("E1").Value = 1.354554295539856
